# Secondary Invoice Template update
# - Adds "Incorrect Product Category" to the Validations lookup list
# - Updates the H2:H983 list validation to include the new row
# - Replaces the sample data row (row 2) on "Sec invoice Master" with a new test record

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Validations sheet: insert "Incorrect Product Category" in alphabetical
#    order right after "Incorrect Freight Class" (old row 6 -> new row 7),
#    pushing the remaining reasons down by one row.
# ---------------------------------------------------------------------------
$wsVal = $wb.Worksheets.Item("Validations")
$wsVal.Range("A7").EntireRow.Insert()
$wsVal.Range("A7").Value = "Incorrect Product Category"

# ---------------------------------------------------------------------------
# 2. Update the "Secondary Category" dropdown (H2:H983) so it still lists the
#    Validations!A1:A12 range (one extra row was added above) and tighten the
#    applied range from H2:H1001 to H2:H983.
# ---------------------------------------------------------------------------
$wsMain = $wb.Worksheets.Item("Sec invoice Master")
$oldValidationRange = $wsMain.Range("H2:H1001")
$oldValidationRange.Validation.Delete()

$newValidationRange = $wsMain.Range("H2:H983")
$newValidationRange.Validation.Add(3, 1, 1, "=Validations!`$A`$1:`$A`$12", 0)
$newValidationRange.Validation.IgnoreBlank = $true
$newValidationRange.Validation.InCellDropdown = $true
$newValidationRange.Validation.ShowInput = $true
$newValidationRange.Validation.ShowError = $true

# ---------------------------------------------------------------------------
# 3. Replace the example data row (row 2) with a new test record.
#    Most of these columns are formatted as text ("@"), so the string values
#    are written verbatim (no date/number auto-conversion). The Customer and
#    COGS columns hold plain numbers, so their formatting is reset to
#    General first so the values are stored as numbers, not text.
# ---------------------------------------------------------------------------
$wsMain.Range("A2").Value = "57410875"              # FC Order ID
$wsMain.Range("B2").Value = "07/07/2021"            # Date CREATED
$wsMain.Range("C2").NumberFormat = "General"
$wsMain.Range("C2").Value = "FCT874922514537512960" # Tracking #
$wsMain.Range("D2").NumberFormat = "General"
$wsMain.Range("D2").Value = "FC Test Carrier"       # Carrier
$wsMain.Range("E2").NumberFormat = "General"
$wsMain.Range("E2").Value = 5679                    # Customer
$wsMain.Range("F2").Value = "200.10"                # New Invoice Amount
$wsMain.Range("G2").NumberFormat = "General"
$wsMain.Range("G2").Value = "Assembly"              # Invoice Reason
$wsMain.Range("H2").NumberFormat = "General"
$wsMain.Range("H2").Value = "Service Upgrade"       # Secondary Category
$wsMain.Range("I2").NumberFormat = "General"
$wsMain.Range("I2").Value = "57410875+12"           # SECONDARY INV #
$wsMain.Range("J2").NumberFormat = "General"
$wsMain.Range("J2").Value = "Yes"                   # On Account?
$wsMain.Range("K2").NumberFormat = "General"
$wsMain.Range("K2").Value = 222                     # COGS

# ---------------------------------------------------------------------------
# 4. Light formatting touch-ups that came along with the re-typed row:
#    the heavy row-wide border is gone except on the Tracking # (C2) and
#    Secondary Category (H2) cells, and the order id / date columns line up
#    differently.
# ---------------------------------------------------------------------------
$wsMain.Range("A2:K2").Borders.LineStyle = -4142    # xlLineStyleNone

$c2 = $wsMain.Range("C2")
$c2.Borders.LineStyle = 1                           # xlContinuous
$c2.Borders.Weight = 2                              # xlThin
$c2.HorizontalAlignment = -4131                     # xlLeft

$h2 = $wsMain.Range("H2")
$h2.VerticalAlignment = -4108                       # xlCenter

$wsMain.Range("B2").HorizontalAlignment = -4152     # xlRight
